$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 46701.855
$ws.Range("J133").Value = 46701.855
$ws.Range("L133").Value = 46701.855
$ws.Range("N133").Value = -56821.855

$ws.Range("H135").Value = 1561.5
$ws.Range("I135").Value = 1561.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 14053.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11518.5
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 125002220
$ws.Range("I137").Value = 142858830
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 428576490
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -428573940
$ws.Range("N137").Value = -23100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3561.375
$ws.Range("I61").Value = 2718.25
$ws.Range("J61").Value = 4404.5
$ws.Range("K61").Value = 2718.25
$ws.Range("L61").Value = 4404.5
$ws.Range("M61").Value = -2506.25
$ws.Range("N61").Value = -4828.5

$ws.Range("H74").Value = 6538.2915
$ws.Range("I74").Value = 1515.5625
$ws.Range("J74").Value = 16583.75
$ws.Range("K74").Value = 1515.5625
$ws.Range("L74").Value = 16583.75
$ws.Range("M74").Value = -641.5625
$ws.Range("N74").Value = -18331.75

$ws.Range("H77").Value = 6538.2915
$ws.Range("I77").Value = 1515.5625
$ws.Range("J77").Value = 16583.75
$ws.Range("K77").Value = 7577.8125
$ws.Range("L77").Value = 82918.75
$ws.Range("M77").Value = -3209.8125
$ws.Range("N77").Value = -91654.75

$ws.Range("H122").Value = 1704.6923
$ws.Range("I122").Value = 1559.6666
$ws.Range("J122").Value = 2031
$ws.Range("K122").Value = 4678.9998
$ws.Range("L122").Value = 6093
$ws.Range("M122").Value = -2228.9998
$ws.Range("N122").Value = -10993

$ws.Range("H132").Value = 2681.6
$ws.Range("I132").Value = 2142.366
$ws.Range("J132").Value = 5138.1113
$ws.Range("K132").Value = 6427.098
$ws.Range("L132").Value = 15414.3339
$ws.Range("M132").Value = -3897.098
$ws.Range("N132").Value = -20474.3339

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3561.375
$ws.Range("I136").Value = 2718.25
$ws.Range("J136").Value = 4404.5
$ws.Range("K136").Value = 8154.75
$ws.Range("L136").Value = 13213.5
$ws.Range("M136").Value = -5604.75
$ws.Range("N136").Value = -18313.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 401368
$ws.Range("I75").Value = 2500
$ws.Range("K75").Value = 2500
$ws.Range("M75").Value = -1564

$ws.Range("H78").Value = 401368
$ws.Range("I78").Value = 2500
$ws.Range("K78").Value = 7500
$ws.Range("M78").Value = -2820

$ws.Range("H107").Value = 2445.889
$ws.Range("I107").Value = 2612.5
$ws.Range("J107").Value = 1113
$ws.Range("K107").Value = 2612.5
$ws.Range("L107").Value = 1113
$ws.Range("M107").Value = -692.5
$ws.Range("N107").Value = -4953

$ws.Range("H132").Value = 44685.715
$ws.Range("J132").Value = 44685.715
$ws.Range("L132").Value = 44685.715
$ws.Range("N132").Value = -54805.715

$ws.Range("H134").Value = 2618.3257
$ws.Range("I134").Value = 1826.3715
$ws.Range("K134").Value = 5479.1145
$ws.Range("M134").Value = -2944.1145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5276.276
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5276.276
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5276.276
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5866.276

$ws.Range("H34").Value = 5276.276
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5276.276
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5276.276
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5680.276

$ws.Range("H58").Value = 52634004
$ws.Range("I58").Value = 76924536
$ws.Range("J58").Value = 4518.5
$ws.Range("K58").Value = 76924536
$ws.Range("L58").Value = 4518.5
$ws.Range("M58").Value = -76924333
$ws.Range("N58").Value = -4924.5

$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 20000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 20000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -18855
$ws.Range("N59").ClearContents()

$ws.Range("H107").Value = 594.3
$ws.Range("I107").Value = 470
$ws.Range("J107").Value = 718.6
$ws.Range("K107").Value = 470
$ws.Range("L107").Value = 718.6
$ws.Range("M107").Value = 1450
$ws.Range("N107").Value = -4558.6

$ws.Range("H122").Value = 1991.4667
$ws.Range("I122").Value = 1246.8823
$ws.Range("J122").Value = 2965.1538
$ws.Range("K122").Value = 3740.6469
$ws.Range("L122").Value = 8895.4614
$ws.Range("M122").Value = -1290.6469
$ws.Range("N122").Value = -13795.4614

$ws.Range("H123").Value = 28548
$ws.Range("J123").Value = 28548
$ws.Range("L123").Value = 28548
$ws.Range("N123").Value = -38348

$ws.Range("H132").Value = 4631427.5
$ws.Range("I132").Value = 5953565.5
$ws.Range("J132").Value = 3943.875
$ws.Range("K132").Value = 17860696.5
$ws.Range("L132").Value = 11831.625
$ws.Range("M132").Value = -17858166.5
$ws.Range("N132").Value = -16891.625

$ws.Range("H134").Value = 33335760
$ws.Range("I134").Value = 38463060
$ws.Range("J134").Value = 26319452
$ws.Range("K134").Value = 115389180
$ws.Range("L134").Value = 78958356
$ws.Range("M134").Value = -115386645
$ws.Range("N134").Value = -78963426

$ws.Range("H136").Value = 52634004
$ws.Range("I136").Value = 76924536
$ws.Range("J136").Value = 4518.5
$ws.Range("K136").Value = 230773608
$ws.Range("L136").Value = 13555.5
$ws.Range("M136").Value = -230771058
$ws.Range("N136").Value = -18655.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5953869.5
$ws.Range("J131").Value = 6946091
$ws.Range("L131").Value = 20838273
$ws.Range("N131").Value = -20848353

$ws.Range("H132").Value = 1751
$ws.Range("I132").Value = 1601.3334
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 14412.0006
$ws.Range("L132").Value = 19800
$ws.Range("M132").Value = -11882.0006
$ws.Range("N132").Value = -24860

$ws.Range("H134").Value = 6016.433
$ws.Range("I134").Value = 3073.3157
$ws.Range("K134").Value = 9219.947100000001
$ws.Range("M134").Value = -4149.947100000001

$ws.Range("H137").Value = 6316112
$ws.Range("I137").Value = 9093702
$ws.Range("J137").Value = 205413.2
$ws.Range("K137").Value = 27281106
$ws.Range("L137").Value = 616239.6000000001
$ws.Range("M137").Value = -27276006
$ws.Range("N137").Value = -626439.6000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2899.394
$ws.Range("I126").Value = 2723.3333
$ws.Range("K126").Value = 8169.999899999999
$ws.Range("M126").Value = -5699.999899999999

$ws.Range("H132").Value = 3054.6597
$ws.Range("I132").Value = 2701.2927
$ws.Range("J132").Value = 5469.3335
$ws.Range("K132").Value = 8103.8781
$ws.Range("L132").Value = 16408.0005
$ws.Range("M132").Value = -5573.8781
$ws.Range("N132").Value = -21468.0005

$ws.Range("H133").Value = 54908.332
$ws.Range("J133").Value = 54908.332
$ws.Range("L133").Value = 54908.332
$ws.Range("N133").Value = -65028.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3552.2856
$ws.Range("I7").Value = 3397.6667
$ws.Range("J7").Value = 3578.0557
$ws.Range("K7").Value = 3397.6667
$ws.Range("L7").Value = 3578.0557
$ws.Range("M7").Value = -3285.6667
$ws.Range("N7").Value = -3802.0557

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H20").Value = 3454.182
$ws.Range("J20").Value = 4665.3335
$ws.Range("L20").Value = 4665.3335
$ws.Range("N20").Value = -5117.3335

$ws.Range("H126").Value = 3552.2856
$ws.Range("I126").Value = 3397.6667
$ws.Range("J126").Value = 3578.0557
$ws.Range("K126").Value = 10193.0001
$ws.Range("L126").Value = 10734.1671
$ws.Range("M126").Value = -7723.000100000001
$ws.Range("N126").Value = -15674.1671

$ws.Range("H136").Value = 8213.647000000001
$ws.Range("I136").Value = 5050.1
$ws.Range("J136").Value = 12733
$ws.Range("K136").Value = 15150.3
$ws.Range("L136").Value = 38199
$ws.Range("M136").Value = -12600.3
$ws.Range("N136").Value = -43299

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2180.25
$ws.Range("I107").Value = 2337.4092
$ws.Range("J107").Value = 451.5
$ws.Range("K107").Value = 7012.2276
$ws.Range("L107").Value = 1354.5
$ws.Range("M107").Value = -5092.2276
$ws.Range("N107").Value = -5194.5

$ws.Range("H132").Value = 4174.9355
$ws.Range("I132").Value = 4903.4443
$ws.Range("J132").Value = 3166.2307
$ws.Range("K132").Value = 14710.3329
$ws.Range("L132").Value = 9498.6921
$ws.Range("M132").Value = -12180.3329
$ws.Range("N132").Value = -14558.6921

$ws.Range("H136").Value = 4010
$ws.Range("I136").Value = 1707.3
$ws.Range("J136").Value = 6312.7
$ws.Range("K136").Value = 5121.9
$ws.Range("L136").Value = 18938.1
$ws.Range("M136").Value = -2571.9
$ws.Range("N136").Value = -24038.1
